# audit-form.docx header fix:
#   - "{{ caseId }" and "}{% if " were split across two separate <w:r> runs
#     (a leftover artefact that even left the placeholder's closing "}}"
#     broken across runs). Re-saving the merged text collapses them back
#     into a single run/""<w:t>"".
#   - the trailing line break after "{% endif %}" lived in its own
#     run; merge it back into the previous run's <w:r> so the <w:br/>
#     is a sibling of the preceding <w:t>, not a separate run.
#
# Both fixes are pure run-merges: the visible text does not change, only
# the run structure the text is stored in. We force Word to re-flow the
# runs by replacing each boundary-spanning bit of text with itself, which
# makes Word collapse adjacent same-formatted runs into one.

$d = $word.ActiveDocument

# --- Fix 1: merge the "...caseId }" / "}{% if ..." runs -------------------
$text = $d.Content.Text
$start = $text.IndexOf("{{ caseId }")
if ($start -lt 0) { throw "Could not find '{{ caseId }' in document" }
$end = $start + "{{ caseId }}{% if ".Length

$rng = $d.Range($start, $end)
if ($rng.Text -ne "{{ caseId }}{% if ") {
    throw "Unexpected text in range: '$($rng.Text)'"
}
$rng.Find.ClearFormatting()
$rng.Find.Execute($rng.Text, $false, $false, $false, $false, $false, $true, 1, $false, $rng.Text, 2) | Out-Null

# --- Fix 2: merge the trailing <w:br/> run back into the previous run -----
$text = $d.Content.Text
$brIdx = $text.IndexOf("{% endif %}") + "{% endif %}".Length
$rng2 = $d.Range($brIdx, $brIdx + 1)
if ([int][char]$rng2.Text[0] -ne 11) {
    throw "Expected a line break character after '{% endif %}', got '$($rng2.Text)'"
}
$joinRng = $d.Range($brIdx - 1, $brIdx + 1)
$joinRng.Find.ClearFormatting()
$joinRng.Find.Execute($joinRng.Text, $false, $false, $false, $false, $false, $true, 1, $false, $joinRng.Text, 2) | Out-Null
